$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AuthTests")
$ws.Activate()

# Update the "RightLoginTestData" table (C24:E30) with a different mix of
# valid/invalid login credentials.
$ws.Range("C25").Value = "admin"
$ws.Range("D25").Value = "admin"

$ws.Range("D27").Value = "invalid"

$ws.Range("C28").Value = "invalid"

$ws.Range("C30").Value = "Admin"
$ws.Range("D30").Value = "admin"

# Leave the cursor where the author last clicked while editing the sheet.
$ws.Range("F26").Select()
